# 99_1.xlsx was refreshed with a re-pulled/re-labeled extract of the same
# confirmations report: row labels become "<Branch>, <Metric>" pairs, the
# stray duplicate "New nominations" row at the bottom is gone (41 -> 40
# rows), a "Total new nominations" row/value is introduced, and the Air
# Force "New nominations" figure becomes a real number instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-removed trailing row so the table is 40 rows (A1:B40).
$ws.Rows(41).Delete()

$ws.Range("A1").Value = "Labels"
$ws.Range("B1").Value = "Values"

$ws.Range("A2").Value = "Congress"
$ws.Range("B2").Value = 99

$ws.Range("A3").Value = "Session"
$ws.Range("B3").Value = 1

$ws.Range("A4").Value = "Start Date"
$ws.Range("B4").Value = 31050

$ws.Range("A5").Value = "End Date"
$ws.Range("B5").Value = 31401

$ws.Range("A6").Value = "Civilian"

$ws.Range("A7").Value = "     Civilian, New nominations"
$ws.Range("B7").Value = 606

$ws.Range("A8").Value = "     Civilian, Confirmed "
$ws.Range("B8").Value = 491

$ws.Range("A9").Value = "     Civilian, Unconfirmed "
$ws.Range("B9").Value = 68

$ws.Range("A10").Value = "     Civilian, Withdrawn"
$ws.Range("B10").Value = 7

$ws.Range("A11").Value = "     Civilian, Failed at Aug adjournment"
$ws.Range("B11").Value = 15

$ws.Range("A12").Value = "     Civilian, Failed at sine die adjournment"
$ws.Range("B12").Value = 19

$ws.Range("A13").Value = "     Civilian, Superseded by recess reappointments"
$ws.Range("B13").Value = 6

$ws.Range("A14").Value = "Civilian (lists)"

$ws.Range("A15").Value = "     Civilian (lists), New nominations"
$ws.Range("B15").Value = 3113

$ws.Range("A16").Value = "     Civilian (lists), Confirmed "
$ws.Range("B16").Value = 3112

$ws.Range("A17").Value = "     Civilian (lists), Unconfirmed "
$ws.Range("B17").Value = 1

$ws.Range("A18").Value = "Air Force "

$ws.Range("A19").Value = "     Air Force, New nominations"
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("B19").Value = 21367

$ws.Range("A20").Value = "     Air Force, Confirmed "
$ws.Range("B20").Value = 19013

$ws.Range("A21").Value = "     Air Force, Unconfirmed "
$ws.Range("B21").Value = 2354

$ws.Range("A22").Value = "Army "

$ws.Range("A23").Value = "     Army, New nominations"
$ws.Range("B23").Value = 15370

$ws.Range("A24").Value = "     Army, Confirmed "
$ws.Range("B24").Value = 14478

$ws.Range("A25").Value = "     Army, Unconfirmed"
$ws.Range("B25").Value = 892

$ws.Range("A26").Value = "Navy "

$ws.Range("A27").Value = "     Navy, New nominations"
$ws.Range("B27").Value = 16721

$ws.Range("A28").Value = "     Navy, Confirmed "
$ws.Range("B28").Value = 16720

$ws.Range("A29").Value = "     Navy, Withdrawn"
$ws.Range("B29").Value = 1

$ws.Range("A30").Value = "Marine Corps "

$ws.Range("A31").Value = "     Marine Corps, New nominations"
$ws.Range("B31").Value = 2466

$ws.Range("A32").Value = "     Marine Corps, Confirmed "
$ws.Range("B32").Value = 2104

$ws.Range("A33").Value = "     Marine Corps, Unconfirmed "
$ws.Range("B33").Value = 362

$ws.Range("A34").Value = "Total new nominations"
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("B34").Value = 59643

$ws.Range("A35").Value = "Total confirmed "
$ws.Range("B35").Value = 55918

$ws.Range("A36").Value = "Total unconfirmed "
$ws.Range("B36").Value = 3677

$ws.Range("A37").Value = "Total withdrawn"
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("B37").Value = 8

$ws.Range("A38").Value = "Total failed at Aug adjournment"
$ws.Range("B38").Value = 15

$ws.Range("A39").Value = "Total failed at sine die adjournment"
$ws.Range("B39").Value = 19

$ws.Range("A40").Value = "Total superseded by recess reappointments"
$ws.Range("B40").Value = 6

$excel.CutCopyMode = $false
